# Add new columns I (I0) and J (IF) to the worksheet, mirroring the
# formatting of the existing header/data columns (A-H).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new headers "I0" and "IF" in I1/J1 ---
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the header style used by the other header cells (bold font,
# thin box border, centered horizontally, top-aligned vertically).
$headerRange = $ws.Range("I1:J1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108   # xlCenter
$headerRange.VerticalAlignment = -4160     # xlTop
$headerRange.Borders.LineStyle = 1         # xlContinuous
$headerRange.Borders.Weight = 2            # xlThin

# --- Data rows (2-77): new numeric values in columns I and J ---
$data = @(
    @(2, 8, 8),
    @(3, 7, 7),
    @(4, 8, 8),
    @(5, 8, 9),
    @(6, 8, 8),
    @(7, 9, 9),
    @(8, 7, 7),
    @(9, 8, 8),
    @(10, 7, 8),
    @(11, 6, 6),
    @(12, 6, 7),
    @(13, 7, 7),
    @(14, 9, 9),
    @(15, 6, 6),
    @(16, 6, 6),
    @(17, 7, 7),
    @(18, 8, 8),
    @(19, 9, 9),
    @(20, 7, 7),
    @(21, 7, 7),
    @(22, 8, 8),
    @(23, 6, 6),
    @(24, 6, 6),
    @(25, 6, 6),
    @(26, 8, 8),
    @(27, 5, 5),
    @(28, 9, 9),
    @(29, 8, 8),
    @(30, 6, 6),
    @(31, 7, 7),
    @(32, 9, 9),
    @(33, 7, 7),
    @(34, 7, 7),
    @(35, 8, 8),
    @(36, 6, 6),
    @(37, 8, 8),
    @(38, 4, 4),
    @(39, 8, 8),
    @(40, 8, 8),
    @(41, 8, 8),
    @(42, 6, 7),
    @(43, 4, 5),
    @(44, 8, 8),
    @(45, 6, 6),
    @(46, 10, 10),
    @(47, 9, 9),
    @(48, 7, 7),
    @(49, 11, 11),
    @(50, 7, 7),
    @(51, 10, 10),
    @(52, 7, 7),
    @(53, 7, 7),
    @(54, 10, 10),
    @(55, 8, 8),
    @(56, 7, 7),
    @(57, 9, 9),
    @(58, 8, 8),
    @(59, 7, 7),
    @(60, 8, 8),
    @(61, 8, 8),
    @(62, 7, 7),
    @(63, 8, 8),
    @(64, 10, 10),
    @(65, 8, 8),
    @(66, 8, 8),
    @(67, 8, 8),
    @(68, 8, 8),
    @(69, 5, 5),
    @(70, 9, 9),
    @(71, 4, 4),
    @(72, 7, 7),
    @(73, 7, 7),
    @(74, 5, 5),
    @(75, 5, 5),
    @(76, 4, 4),
    @(77, 8, 8),
)

foreach ($row in $data) {
    $r = $row[0]
    $iVal = $row[1]
    $jVal = $row[2]
    $ws.Cells.Item($r, 9).Value = $iVal
    $ws.Cells.Item($r, 10).Value = $jVal
}

